# Scheduled-runner market data refresh: updates currentAveragePrice / Leve
# price / profit columns (H:N) for a handful of rows across several
# crafter-job sheets. Values are plain literals (no formulas in this sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K15").Value = 7977.599999999999
$ws.Range("I15").Value = 2659.2
$ws.Range("M15").Value = -7808.599999999999
$ws.Range("H15").Value = 2659.2
$ws.Range("J107").Value = 166668670
$ws.Range("N107").Value = -166672510
$ws.Range("K107").Value = 5683788.5
$ws.Range("I107").Value = 5683788.5
$ws.Range("M107").Value = -5681868.5
$ws.Range("H107").Value = 25001974
$ws.Range("L107").Value = 166668670
$ws.Range("J113").Value = 100017140
$ws.Range("N113").Value = -100023648
$ws.Range("H113").Value = 125013384
$ws.Range("L113").Value = 100017140
$ws.Range("J138").Value = 5354.479
$ws.Range("N138").Value = -26343.437
$ws.Range("K138").Value = 10472.0625
$ws.Range("I138").Value = 3490.6875
$ws.Range("M138").Value = -5332.0625
$ws.Range("H138").Value = 5011.7124
$ws.Range("L138").Value = 16063.437

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J2").Value = 500002500
$ws.Range("N2").Value = -500002726
$ws.Range("K2").Value = 2422.6
$ws.Range("I2").Value = 2422.6
$ws.Range("M2").Value = -2309.6
$ws.Range("H2").Value = 58825960
$ws.Range("L2").Value = 500002500
$ws.Range("J74").Value = 5564.923
$ws.Range("N74").Value = -7312.923
$ws.Range("K74").Value = 31805.94
$ws.Range("I74").Value = 31805.94
$ws.Range("M74").Value = -30931.94
$ws.Range("H74").Value = 24390
$ws.Range("L74").Value = 5564.923
$ws.Range("J77").Value = 5564.923
$ws.Range("N77").Value = -36560.615
$ws.Range("K77").Value = 159029.7
$ws.Range("I77").Value = 31805.94
$ws.Range("M77").Value = -154661.7
$ws.Range("H77").Value = 24390
$ws.Range("L77").Value = 27824.615
$ws.Range("J110").Value = 111111784
$ws.Range("N110").Value = -111115874
$ws.Range("K110").Value = 16542.857
$ws.Range("I110").Value = 16542.857
$ws.Range("M110").Value = -14497.857
$ws.Range("H110").Value = 33345114
$ws.Range("L110").Value = 111111784
$ws.Range("J116").Value = 500002500
$ws.Range("N116").Value = -500007088
$ws.Range("K116").Value = 2422.6
$ws.Range("I116").Value = 2422.6
$ws.Range("M116").Value = -128.5999999999999
$ws.Range("H116").Value = 58825960
$ws.Range("L116").Value = 500002500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 500002500
$ws.Range("N3").Value = -500002728
$ws.Range("K3").Value = 2422.6
$ws.Range("I3").Value = 2422.6
$ws.Range("M3").Value = -2308.6
$ws.Range("H3").Value = 58825960
$ws.Range("L3").Value = 500002500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K31").Value = 2069.5789
$ws.Range("I31").Value = 2069.5789
$ws.Range("M31").Value = -1774.5789
$ws.Range("H31").Value = 5052.6777
$ws.Range("K34").Value = 2069.5789
$ws.Range("I34").Value = 2069.5789
$ws.Range("M34").Value = -1867.5789
$ws.Range("H34").Value = 5052.6777
$ws.Range("J107").Value = 3450.3333
$ws.Range("N107").Value = -7290.3333
$ws.Range("H107").Value = 2494.9033
$ws.Range("L107").Value = 3450.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J11").Value = 4560
$ws.Range("N11").Value = -13960
$ws.Range("H11").Value = 2868.375
$ws.Range("L11").Value = 13680
$ws.Range("J82").Value = 100000
$ws.Range("N82").Value = -300812
$ws.Range("H82").Value = 61004
$ws.Range("L82").Value = 300000
$ws.Range("J85").Value = 100000
$ws.Range("N85").Value = -302808
$ws.Range("H85").Value = 61004
$ws.Range("L85").Value = 300000
$ws.Range("J107").Value = 4884467.5
$ws.Range("N107").Value = -14657242.5
$ws.Range("K107").Value = 6001087.199999999
$ws.Range("I107").Value = 2000362.4
$ws.Range("M107").Value = -5999167.199999999
$ws.Range("H107").Value = 4584040
$ws.Range("L107").Value = 14653402.5
$ws.Range("K114").Value = 1225.00002
$ws.Range("I114").Value = 408.33334
$ws.Range("M114").Value = 2028.99998
$ws.Range("H114").Value = 1746
$ws.Range("J137").Value = 999999
$ws.Range("N137").Value = -3010197
$ws.Range("H137").Value = 289722.28
$ws.Range("L137").Value = 2999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K102").Value = 3674.8076
$ws.Range("I102").Value = 3674.8076
$ws.Range("M102").Value = -2052.8076
$ws.Range("H102").Value = 3874.6072
$ws.Range("K122").Value = 24149739
$ws.Range("I122").Value = 8049913
$ws.Range("M122").Value = -24147289
$ws.Range("H122").Value = 8049913
$ws.Range("J123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("K126").Value = 14997
$ws.Range("I126").Value = 4999
$ws.Range("M126").Value = -12527
$ws.Range("H126").Value = 9461.299999999999
$ws.Range("J141").Value = 70000
$ws.Range("N141").Value = -80360
$ws.Range("H141").Value = 70000
$ws.Range("L141").Value = 70000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 11111
$ws.Range("N22").Value = -11701
$ws.Range("H22").Value = 8555.5
$ws.Range("L22").Value = 11111
$ws.Range("J27").Value = 11111
$ws.Range("N27").Value = -11325
$ws.Range("H27").Value = 8555.5
$ws.Range("L27").Value = 11111
$ws.Range("J40").Value = 8642.714
$ws.Range("N40").Value = -8914.714
$ws.Range("K40").Value = 4739.8
$ws.Range("I40").Value = 4739.8
$ws.Range("M40").Value = -4603.8
$ws.Range("H40").Value = 7016.5
$ws.Range("L40").Value = 8642.714
$ws.Range("J127").Value = 74810
$ws.Range("N127").Value = -84730
$ws.Range("H127").Value = 74810
$ws.Range("L127").Value = 74810
$ws.Range("K136").Value = 11230.5
$ws.Range("I136").Value = 3743.5
$ws.Range("M136").Value = -8680.5
$ws.Range("H136").Value = 10815.259

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J136").Value = 2462.2693
$ws.Range("N136").Value = -12486.8079
$ws.Range("K136").Value = 5556.9729
$ws.Range("I136").Value = 1852.3243
$ws.Range("M136").Value = -3006.9729
$ws.Range("H136").Value = 2010.91
$ws.Range("L136").Value = 7386.8079
$ws.Range("J140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()
$ws.Range("K141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("H141").Value = 0
$ws.Range("L141").Value = 0
